$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 989
$ws.Range("I28").Value = 850
$ws.Range("J28").Value = 1406
$ws.Range("K28").Value = 850
$ws.Range("L28").Value = 1406
$ws.Range("M28").Value = -365
$ws.Range("N28").Value = -2376

# Sheet ALC, Row 37
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H37").Value = 196.66667
$ws.Range("I37").Value = 140
$ws.Range("J37").Value = 225
$ws.Range("K37").Value = 420
$ws.Range("L37").Value = 675
$ws.Range("M37").Value = -294
$ws.Range("N37").Value = -927

# Sheet ALC, Row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3381.0278
$ws.Range("I64").Value = 3259.8635
$ws.Range("J64").Value = 3571.4285
$ws.Range("K64").Value = 3259.8635
$ws.Range("L64").Value = 3571.4285
$ws.Range("M64").Value = -3011.8635
$ws.Range("N64").Value = -4067.4285

# Sheet ALC, Row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3381.0278
$ws.Range("I67").Value = 3259.8635
$ws.Range("J67").Value = 3571.4285
$ws.Range("K67").Value = 3259.8635
$ws.Range("L67").Value = 3571.4285
$ws.Range("M67").Value = -2401.8635
$ws.Range("N67").Value = -5287.4285

# Sheet ALC, Row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1845.6364
$ws.Range("J112").Value = 1950.2
$ws.Range("L112").Value = 5850.6
$ws.Range("N112").Value = -8066.6

# Sheet ALC, Row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 40002740
$ws.Range("J116").Value = 3566.6667
$ws.Range("L116").Value = 3566.6667
$ws.Range("N116").Value = -10450.6667

# Sheet ALC, Row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1060.8372
$ws.Range("I129").Value = 726.8570999999999
$ws.Range("J129").Value = 1125.7778
$ws.Range("K129").Value = 2180.5713
$ws.Range("L129").Value = 3377.3334
$ws.Range("M129").Value = 2819.4287
$ws.Range("N129").Value = -13377.3334

# Sheet ALC, Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1433.1515
$ws.Range("I137").Value = 1291.7307
$ws.Range("J137").Value = 1958.4286
$ws.Range("K137").Value = 3875.1921
$ws.Range("L137").Value = 5875.2858
$ws.Range("M137").Value = -1325.1921
$ws.Range("N137").Value = -10975.2858

# Sheet ALC, Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2535210.8
$ws.Range("I138").Value = 9092422
$ws.Range("J138").Value = 4357.456
$ws.Range("K138").Value = 27277266
$ws.Range("L138").Value = 13072.368
$ws.Range("M138").Value = -27272126
$ws.Range("N138").Value = -23352.368

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1624.5
$ws.Range("I61").Value = 1624.5
$ws.Range("K61").Value = 1624.5
$ws.Range("M61").Value = -1412.5

# Sheet ARM, Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1454.2
$ws.Range("I122").Value = 1454.2
$ws.Range("K122").Value = 4362.6
$ws.Range("M122").Value = -1912.6

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2015.4286
$ws.Range("I132").Value = 1584.5
$ws.Range("K132").Value = 4753.5
$ws.Range("M132").Value = -2223.5

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1624.5
$ws.Range("I136").Value = 1624.5
$ws.Range("K136").Value = 4873.5
$ws.Range("M136").Value = -2323.5

# Sheet BSM, Row 53
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 24490
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# Sheet CRP, Row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -9264
$ws.Range("N51").ClearContents()

# Sheet CRP, Row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 31567.857
$ws.Range("I59").Value = 10200
$ws.Range("K59").Value = 10200
$ws.Range("M59").Value = -9055

# Sheet CRP, Row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -9652
$ws.Range("N61").ClearContents()

# Sheet CRP, Row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 16586
$ws.Range("I69").Value = 16586
$ws.Range("K69").Value = 16586
$ws.Range("M69").Value = -15837

# Sheet CRP, Row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 16586
$ws.Range("I72").Value = 16586
$ws.Range("K72").Value = 49758
$ws.Range("M72").Value = -46014

# Sheet CRP, Row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 447.73334
$ws.Range("I107").Value = 362.75
$ws.Range("K107").Value = 362.75
$ws.Range("M107").Value = 1557.25

# Sheet CRP, Row 137
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 29195
$ws.Range("I137").Value = 29000
$ws.Range("J137").Value = 29780
$ws.Range("K137").Value = 29000
$ws.Range("L137").Value = 29780
$ws.Range("M137").Value = -23900
$ws.Range("N137").Value = -39980

# Sheet CUL, Row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1207.7142
$ws.Range("I5").Value = 1473.1
$ws.Range("K5").Value = 4419.299999999999
$ws.Range("M5").Value = -4307.299999999999

# Sheet CUL, Row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 162.375
$ws.Range("I26").Value = 162.375
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 487.125
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -199.125
$ws.Range("N26").ClearContents()

# Sheet CUL, Row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 4285.7144
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4285.7144
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 12857.1432
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -13113.1432

# Sheet CUL, Row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 661.2069
$ws.Range("I113").Value = 645.1818
$ws.Range("J113").Value = 671
$ws.Range("K113").Value = 1935.5454
$ws.Range("L113").Value = 2013
$ws.Range("M113").Value = 234.4546
$ws.Range("N113").Value = -6353

# Sheet CUL, Row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 961.7692
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 991.1818
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 2973.5454
$ws.Range("M121").Value = -1090
$ws.Range("N121").Value = -5593.5454

# Sheet CUL, Row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2746.25
$ws.Range("I125").Value = 1860
$ws.Range("J125").Value = 3149.0908
$ws.Range("K125").Value = 5580
$ws.Range("L125").Value = 9447.2724
$ws.Range("M125").Value = -660
$ws.Range("N125").Value = -19287.2724

# Sheet CUL, Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17860284
$ws.Range("I131").Value = 9434.546
$ws.Range("J131").Value = 22223824
$ws.Range("K131").Value = 28303.638
$ws.Range("L131").Value = 66671472
$ws.Range("M131").Value = -23263.638
$ws.Range("N131").Value = -66681552

# Sheet CUL, Row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1207.7142
$ws.Range("I135").Value = 1473.1
$ws.Range("K135").Value = 13257.9
$ws.Range("M135").Value = -10722.9

# Sheet GSM, Row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I4").Value = 20000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -19888
$ws.Range("N4").ClearContents()

# Sheet LTW, Row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2685.25
$ws.Range("I136").Value = 2367.28
$ws.Range("J136").Value = 5335
$ws.Range("K136").Value = 7101.84
$ws.Range("L136").Value = 16005
$ws.Range("M136").Value = -4551.84
$ws.Range("N136").Value = -21105

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 952.3333
$ws.Range("I136").Value = 930.5862
$ws.Range("J136").Value = 1110
$ws.Range("K136").Value = 2791.7586
$ws.Range("L136").Value = 3330
$ws.Range("M136").Value = -241.7586000000001
$ws.Range("N136").Value = -8430
